# Auto-generated edit script applying crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Cell = 'D2'; Value = '66.670.83'; ForceText = $false }
    @{ Cell = 'E2'; Value = '  -1.13%  '; ForceText = $false }
    @{ Cell = 'D3'; Value = '3.515.30'; ForceText = $false }
    @{ Cell = 'E3'; Value = '  +0.57%  '; ForceText = $false }
    @{ Cell = 'E4'; Value = '  +0.00%  '; ForceText = $false }
    @{ Cell = 'D5'; Value = '584.25'; ForceText = $true }
    @{ Cell = 'E5'; Value = '  -2.19%  '; ForceText = $false }
    @{ Cell = 'D6'; Value = '175.69'; ForceText = $true }
    @{ Cell = 'E6'; Value = '  -2.64%  '; ForceText = $false }
    @{ Cell = 'E7'; Value = '  +0.03%  '; ForceText = $false }
    @{ Cell = 'D8'; Value = '3.509.13'; ForceText = $false }
    @{ Cell = 'E8'; Value = '  +0.38%  '; ForceText = $false }
    @{ Cell = 'D9'; Value = '0.596'; ForceText = $true }
    @{ Cell = 'E9'; Value = '  -1.89%  '; ForceText = $false }
    @{ Cell = 'D10'; Value = '0.135'; ForceText = $true }
    @{ Cell = 'E10'; Value = '  -2.39%  '; ForceText = $false }
    @{ Cell = 'D11'; Value = '6.91'; ForceText = $true }
    @{ Cell = 'E11'; Value = '  -2.03%  '; ForceText = $false }
    @{ Cell = 'E12'; Value = '  -3.15%  '; ForceText = $false }
    @{ Cell = 'D13'; Value = '4.110.75'; ForceText = $false }
    @{ Cell = 'E13'; Value = '  +0.42%  '; ForceText = $false }
    @{ Cell = 'D14'; Value = '30.49'; ForceText = $true }
    @{ Cell = 'E14'; Value = '  -5.70%  '; ForceText = $false }
    @{ Cell = 'E15'; Value = '  -1.69%  '; ForceText = $false }
    @{ Cell = 'D16'; Value = '66.635.70'; ForceText = $false }
    @{ Cell = 'E16'; Value = '  -1.14%  '; ForceText = $false }
    @{ Cell = 'D17'; Value = '0.0000174'; ForceText = $true }
    @{ Cell = 'E17'; Value = '  -2.39%  '; ForceText = $false }
    @{ Cell = 'D18'; Value = '3.506.18'; ForceText = $false }
    @{ Cell = 'E18'; Value = '  +0.37%  '; ForceText = $false }
    @{ Cell = 'D19'; Value = '6.06'; ForceText = $true }
    @{ Cell = 'E19'; Value = '  -3.97%  '; ForceText = $false }
    @{ Cell = 'D20'; Value = '13.99'; ForceText = $true }
    @{ Cell = 'E20'; Value = '  -2.18%  '; ForceText = $false }
    @{ Cell = 'D21'; Value = '381.43'; ForceText = $true }
    @{ Cell = 'E21'; Value = '  -2.31%  '; ForceText = $false }
    @{ Cell = 'D22'; Value = '7.89'; ForceText = $true }
    @{ Cell = 'E22'; Value = '  -0.84%  '; ForceText = $false }
    @{ Cell = 'D23'; Value = '0.548'; ForceText = $true }
    @{ Cell = 'E23'; Value = '  +0.99%  '; ForceText = $false }
    @{ Cell = 'D24'; Value = '1.00'; ForceText = $true }
    @{ Cell = 'E24'; Value = '  +0.26%  '; ForceText = $false }
    @{ Cell = 'D25'; Value = '72.33'; ForceText = $true }
    @{ Cell = 'E25'; Value = '  -2.23%  '; ForceText = $false }
    @{ Cell = 'D26'; Value = '5.75'; ForceText = $true }
    @{ Cell = 'E26'; Value = '  +0.21%  '; ForceText = $false }
    @{ Cell = 'D27'; Value = '0.0000121'; ForceText = $true }
    @{ Cell = 'E27'; Value = '  -1.04%  '; ForceText = $false }
    @{ Cell = 'D28'; Value = '9.88'; ForceText = $true }
    @{ Cell = 'E28'; Value = '  -4.97%  '; ForceText = $false }
    @{ Cell = 'D29'; Value = '0.173'; ForceText = $true }
    @{ Cell = 'E29'; Value = '  -2.04%  '; ForceText = $false }
    @{ Cell = 'E30'; Value = '  +0.09%  '; ForceText = $false }
    @{ Cell = 'D31'; Value = '24.58'; ForceText = $true }
    @{ Cell = 'E31'; Value = '  +4.26%  '; ForceText = $false }
    @{ Cell = 'D32'; Value = '5.91'; ForceText = $true }
    @{ Cell = 'E32'; Value = '  -4.45%  '; ForceText = $false }
    @{ Cell = 'D33'; Value = '2.02'; ForceText = $true }
    @{ Cell = 'E33'; Value = '  -2.69%  '; ForceText = $false }
    @{ Cell = 'D34'; Value = '1.35'; ForceText = $true }
    @{ Cell = 'E34'; Value = '  -5.23%  '; ForceText = $false }
    @{ Cell = 'E35'; Value = '  -0.03%  '; ForceText = $false }
    @{ Cell = 'D36'; Value = '7.26'; ForceText = $true }
    @{ Cell = 'E36'; Value = '  -1.89%  '; ForceText = $false }
    @{ Cell = 'E37'; Value = '  -1.60%  '; ForceText = $false }
    @{ Cell = 'B38'; Value = 'EnergySwap'; ForceText = $false }
    @{ Cell = 'C38'; Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; ForceText = $false }
    @{ Cell = 'D38'; Value = '30.09'; ForceText = $true }
    @{ Cell = 'E38'; Value = '  +13.76%  '; ForceText = $false }
    @{ Cell = 'B39'; Value = 'Monero'; ForceText = $false }
    @{ Cell = 'C39'; Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; ForceText = $false }
    @{ Cell = 'D39'; Value = '161.50'; ForceText = $true }
    @{ Cell = 'E39'; Value = '  -1.15%  '; ForceText = $false }
    @{ Cell = 'D40'; Value = '0.898'; ForceText = $true }
    @{ Cell = 'E40'; Value = '  +3.10%  '; ForceText = $false }
    @{ Cell = 'D41'; Value = '1.79'; ForceText = $true }
    @{ Cell = 'E41'; Value = '  -5.09%  '; ForceText = $false }
    @{ Cell = 'E42'; Value = '  -4.54%  '; ForceText = $false }
    @{ Cell = 'D43'; Value = '4.50'; ForceText = $true }
    @{ Cell = 'E43'; Value = '  -3.15%  '; ForceText = $false }
    @{ Cell = 'D44'; Value = '2.729.56'; ForceText = $false }
    @{ Cell = 'E44'; Value = '  -4.35%  '; ForceText = $false }
    @{ Cell = 'D45'; Value = '2.54'; ForceText = $true }
    @{ Cell = 'E45'; Value = '  -9.93%  '; ForceText = $false }
    @{ Cell = 'D46'; Value = '0.0703'; ForceText = $true }
    @{ Cell = 'E46'; Value = '  -2.74%  '; ForceText = $false }
    @{ Cell = 'D47'; Value = '40.73'; ForceText = $true }
    @{ Cell = 'E47'; Value = '  -2.32%  '; ForceText = $false }
    @{ Cell = 'D48'; Value = '25.03'; ForceText = $true }
    @{ Cell = 'E48'; Value = '  -6.86%  '; ForceText = $false }
    @{ Cell = 'D49'; Value = '0.0293'; ForceText = $true }
    @{ Cell = 'E49'; Value = '  -2.82%  '; ForceText = $false }
    @{ Cell = 'D50'; Value = '324.52'; ForceText = $true }
    @{ Cell = 'E50'; Value = '  -2.73%  '; ForceText = $false }
    @{ Cell = 'D51'; Value = '1.02'; ForceText = $true }
    @{ Cell = 'E51'; Value = '  -4.03%  '; ForceText = $false }
)

foreach ($chg in $changes) {
    $rng = $ws.Range($chg.Cell)
    if ($chg.ForceText) {
        $rng.NumberFormat = "@"
    }
    $rng.Value = $chg.Value
}
